# Handles float input without breaking stuff
#
# The underlying question-count / scoring computation previously produced
# bad (e.g. float-derived/duplicated) results; this regenerates the
# "quiz" marksheet with the corrected values:
#   - 28 real questions (not 56)
#   - Right/Wrong/NotAttempt = 11/4/13, Max = 28
#   - Marking scheme +4/-1/0, Total = 44/-4, Score = 40/112
#   - the redundant 3rd "Student Ans/Correct Ans" block (columns G:H) removed
#   - the leftover columns D:E trimmed down to just the 3 extra questions
#     that don't fit in the A:B block (rows 16-18), the rest cleared
#   - A column (Student Ans) filled in / re-graded with correct/incorrect
#     styling for every answered question

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# Summary block (rows 10-12)
# ---------------------------------------------------------------------
$ws.Range("A10").Value = "No."
$ws.Range("A10").Style = "mtitleStyle"
$ws.Range("B10").Value = 11
$ws.Range("C10").Value = 4
$ws.Range("D10").Value = 13
$ws.Range("E10").Value = 28

$ws.Range("A11").Value = "Marking"
$ws.Range("A11").Style = "mtitleStyle"
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -1

$ws.Range("A12").Value = "Total"
$ws.Range("A12").Style = "mtitleStyle"
$ws.Range("B12").Value = 44
$ws.Range("C12").Value = -4
$ws.Range("E12").Value = "40/112"

# ---------------------------------------------------------------------
# Drop the third (now unused) Student Ans / Correct Ans block entirely.
# ---------------------------------------------------------------------
$ws.Columns("G:H").Delete()

# ---------------------------------------------------------------------
# Rows 16-18: three extra questions (26-28) kept in columns D:E, with
# the Student Ans side (D) filled in and graded.
# ---------------------------------------------------------------------
$ws.Range("D16").Value = "Option A"
$ws.Range("D16").Style = "correctStyle"

$ws.Range("D17").Value = "Option C"
$ws.Range("D17").Style = "correctStyle"

$ws.Range("D18").Value = "Option A"
$ws.Range("D18").Style = "incorrectStyle"

# ---------------------------------------------------------------------
# Column A (Student Ans) for questions 1-25, graded against column B.
# ---------------------------------------------------------------------
$ws.Range("A18").Value = "Option B"
$ws.Range("A18").Style = "correctStyle"

$ws.Range("A19").Value = "Option C"
$ws.Range("A19").Style = "correctStyle"

$ws.Range("A23").Value = "Option C"
$ws.Range("A23").Style = "incorrectStyle"

$ws.Range("A26").Value = "Option D"
$ws.Range("A26").Style = "incorrectStyle"

$ws.Range("A27").Value = "Option A"
$ws.Range("A27").Style = "correctStyle"

$ws.Range("A30").Value = "Option B"
$ws.Range("A30").Style = "correctStyle"

$ws.Range("A32").Value = "Option C"
$ws.Range("A32").Style = "correctStyle"

$ws.Range("A33").Value = "Option D"
$ws.Range("A33").Style = "correctStyle"

$ws.Range("A34").Value = "Option A"
$ws.Range("A34").Style = "incorrectStyle"

$ws.Range("A35").Value = "Option D"
$ws.Range("A35").Style = "correctStyle"

$ws.Range("A38").Value = "Option A"
$ws.Range("A38").Style = "correctStyle"

$ws.Range("A39").Value = "Option D"
$ws.Range("A39").Style = "correctStyle"

# ---------------------------------------------------------------------
# Everything left over in D19:E40 was only relevant to the removed
# 56-question dataset - clear it out (cells disappear entirely).
# ---------------------------------------------------------------------
$ws.Range("D19:E40").Clear()
